$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows (old row 2 "age" ... row 19 "rm") down by one
$ws.Rows.Item(2).Insert()

# Put the new header value "n" into the freshly inserted row's A cell
$ws.Range("A2").Value = "n"

# Match the cursor position recorded in the saved file
$ws.Range("B2").Select()
